# Insert a new data row before current row 75, shifting existing rows 75-104
# down to 76-105, and populate the new row 75 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75 (pushes old row 75 and below down by one).
# The inserted row inherits the formatting (incl. the date number format in
# column D) from the row above it, same as a manual row insert in Excel.
$ws.Rows(75).Insert()

# Fill in the new row's data
$ws.Range("A75").Value = 6
$ws.Range("B75").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C75").Value = "Metropolitana"
$ws.Range("D75").Value = 45141
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 100112035
$ws.Range("G75").Value = "Bruselas (repollito)"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 220
$ws.Range("K75").Value = 17000
$ws.Range("L75").Value = 18000
$ws.Range("M75").Value = 17545
$ws.Range("N75").Value = "`$/malla 15 kilos"
$ws.Range("O75").Value = "Provincia de Quillota"
$ws.Range("P75").Value = 1170
$ws.Range("Q75").Value = 15
$ws.Range("R75").Value = "Hortaliza"
